$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "practice_areas" sheet changes  (do this first; "overview" must be the
#    last sheet touched/selected so it remains the tab-selected sheet)
# ---------------------------------------------------------------------------
$pa = $wb.Worksheets.Item("practice_areas")

# -- data rows --------------------------------------------------------------
# New shared strings must be appended in the same order they first appear in
# the target workbook (indices 78-87): C5,C6,C7,C8,C9,C2,C3,C4,F1,G1

# Row 5
$pa.Range("A5").Value = 5.1
$pa.Range("B5").Value = 5
$pa.Range("C5").Value = "Contract mobilisation and transition"
$pa.Range("E5").Value = 11
$pa.Range("F5").Value = 0

# Row 6
$pa.Range("A6").Value = 5.2
$pa.Range("B6").Value = 5
$pa.Range("C6").Value = "Managing delivery and performance"
$pa.Range("E6").Value = 14
$pa.Range("F6").Value = 0

# Row 7
$pa.Range("A7").Value = 5.3
$pa.Range("B7").Value = 5
$pa.Range("C7").Value = "Managing contract risk"
$pa.Range("E7").Value = 13
$pa.Range("F7").Value = 0

# Row 8
$pa.Range("A8").Value = 5.4
$pa.Range("B8").Value = 5
$pa.Range("C8").Value = "Supply chain, inventory management/stock control (NHS Only)"
$pa.Range("E8").Value = 12
$pa.Range("F8").Value = 0

# Row 9
$pa.Range("A9").Value = 5.5
$pa.Range("B9").Value = 5
$pa.Range("C9").Value = "Supply chain, logistics (NHS Only)"
$pa.Range("E9").Value = 6
$pa.Range("F9").Value = 0

# Row 2
$pa.Range("A2").Value = 4.1
$pa.Range("B2").Value = 4
$pa.Range("C2").Value = "Applying effective contract terms"
$pa.Range("E2").Value = 10
$pa.Range("F2").Value = 5

# Row 3
$pa.Range("A3").Value = 4.2
$pa.Range("B3").Value = 4
$pa.Range("C3").Value = "Procurement/commercial activity, bid evaluation and supplier selection"
$pa.Range("E3").Value = 10
$pa.Range("F3").Value = 10

# Row 4
$pa.Range("A4").Value = 4.3
$pa.Range("B4").Value = 4
$pa.Range("C4").Value = "Appropriate risk allocation between parties"
$pa.Range("E4").Value = 6
$pa.Range("F4").Value = 0

# -- header row (new columns) ------------------------------------------------
$pa.Range("F1").Value = "pa_complete"
$pa.Range("G1").Value = "percent_complete"

# -- formulas for the new percent_complete column ----------------------------
$pa.Range("G2").Formula = "=(F2/E2)*100"
$pa.Range("G3").Formula = "=(F3/E3)*100"
$pa.Range("G4").Formula = "=(F4/E4)*100"
$pa.Range("G5").Formula = "=(F5/E5)*100"
$pa.Range("G6").Formula = "=(F6/E6)*100"
$pa.Range("G7").Formula = "=(F7/E7)*100"
$pa.Range("G8").Formula = "=(F8/E8)*100"
$pa.Range("G9").Formula = "=(F9/E9)*100"

# -- cell formatting (style index 1 = default font, horizontal left) --------
$pa.Range("A1:G1").HorizontalAlignment = -4131
$pa.Range("A2:C9").HorizontalAlignment = -4131
$pa.Range("E2:F9").HorizontalAlignment = -4131

# -- number format + alignment (style index 5) for the percent column -------
$pa.Range("G2:G9").HorizontalAlignment = -4131
$pa.Range("G2:G9").NumberFormat = "0"

# -- column widths -----------------------------------------------------------
# (target widths 14.42578125 / 51.5703125 / 20.7109375 are not exactly
#  representable by this runtime's column-width engine, which snaps to the
#  nearest 1/6 character; the closest achievable width is used instead)
$pa.Columns.Item(1).ColumnWidth = 13.67
$pa.Columns.Item(2).ColumnWidth = 13.67
$pa.Columns.Item(4).ColumnWidth = 50.7
$pa.Columns.Item(7).ColumnWidth = 19.8

$pa.Range("F4").Select()

# ---------------------------------------------------------------------------
# 2. "overview" sheet changes (done last so it stays the active/tab-selected
#    sheet, matching the target workbook)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("overview")
$overview.Range("E6").Value = 0
$overview.Range("F6").Value = "NOT STARTED"
$overview.Range("H17").Select()
